$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    4  = @{ C = 0.04039603960396039;  E = 0.08992499249924991;  F = -0.07599159915991598; H = -0.06741074107410741; I = -0.03313531353135313; J = 0.0547542332253334 }
    5  = @{ C = -0.1364656465646565;  E = -0.104026402640264;   F = 0.04314431443144313;  H = -0.06672667266726671; I = -0.1426222622262226;  J = -0.1641245224361294 }
    6  = @{ C = 0.05422142214221421; E = -0.08795679567956795; F = 0.1042424242424242;   H = -0.0144974497449745;  I = -0.04927692769276927; J = 0.09220026523032607 }
    7  = @{ C = -0.00468046804680468; E = -0.1267446744674467;  F = 0.06669066906690668;  H = 0.006540654065406539; I = 0.007704770477047704;  J = 0.0592960590228265 }
    8  = @{ C = 0.0396039603960396;  E = -0.224062406240624;   F = 0.5241284128412841;   H = 0.9999999999999999;  I = 0.8634503450345034;   J = 0.01390783957830226 }
    9  = @{ C = 0.9734173417341733;  E = 0.2249264926492649;   F = -0.1854545454545454;  H = -0.0007320732073207321; I = 0.09438943894389439; J = 0.05583562031997463 }
    10 = @{ C = -0.08555655565556554; E = -0.2532253225322532;  F = 0.1829102910291029;   H = 0.05572157215721572;  I = -0.03294329432943294; J = -0.1850073087748675 }
    11 = @{ C = -0.05341734173417341; E = 0.111047104710471;    F = -0.1604800480048005;  H = -0.2049324932493249;  I = -0.1291929192919292;  J = 0.1095505203932362 }
    12 = @{ C = 0.005496549654965496; E = -0.01933393339333933; F = -0.1438703870387039;  H = -0.09008100810081007; I = -0.1041584158415841;  J = -0.127159106917644 }
    13 = @{ C = 0.2603540354035404;  E = 0.6885808580858086;   F = -0.8300870087008699;  H = -0.02807080708070807; I = 0.4652265226522652;   J = 0.1520970948612866 }
    14 = @{ C = -0.1017701770177018; E = -0.1318331833183318;  F = 0.07571557155715571;  H = 0.03169516951695169;  I = -0.03866786678667866; J = -0.2162533881039172 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
